$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.449.72"
$ws.Range("E2").Value = "  +4.33%  "
$ws.Range("D3").Value = "1.592.31"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.66"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.03"
$ws.Range("E8").Value = "  +8.56%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0601"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("D12").Value = "1.818.06"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "1.591.11"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("D16").Value = "28.457.33"
$ws.Range("E16").Value = "  +4.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.21"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.12"
$ws.Range("E18").Value = "  +7.29%  "
$ws.Range("D19").Value = "0.0₃0710"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.45"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.25"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.32"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.62"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0475"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "1.420.82"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.53"
$ws.Range("E39").Value = "  +8.20%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.545"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.84"
$ws.Range("E44").Value = "  +6.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.977"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.79"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "1.730.64"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.72"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").Value = "0.0₆0107"
$ws.Range("E50").Value = "  +6.03%  "
$ws.Range("E51").Value = "  -0.43%  "
